{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the empty paragraph that sits right after \"Just need to now: - \"\n// and directly before the \"Figure out why delete button\" list item, and\n// remove it.\nfor (let i = 0; i < items.length - 1; i++) {\n  const cur = items[i];\n  const next = items[i + 1];\n  if (cur.text.trim() === \"\" && next.text.indexOf(\"Figure out why delete button\") !== -1) {\n    cur.delete();\n    break;\n  }\n}\n\n// Replace the three runs of the \"Figure out why delete button...\" item\n// with a single sentence describing the new state of the bug.\nfor (let i = 0; i < items.length; i++) {\n  const p = items[i];\n  if (p.text.indexOf(\"Figure out why delete button\") !== -1) {\n    p.insertText(\n      \"Delete button is now working for first to do list item but not for later items.\",\n      Word.InsertLocation.replace\n    );\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the empty paragraph that sits between \"Just need to now: - \" and\n# the \"Figure out why delete button\" list item.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $txt = $para.Range.Text.Trim()\n    if ($txt -eq \"\" -and $i -lt $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($i + 1).Range.Text\n        if ($nextText -like \"*Figure out why delete button*\") {\n            $para.Range.Delete()\n            break\n        }\n    }\n}\n\n# Replace the whole \"Figure out why delete button...\" list item text with\n# the updated status sentence.\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"*Figure out why delete button*\") {\n        $r = $para.Range\n        $r.MoveEnd(1, -1) | Out-Null\n        $r.Text = \"Delete button is now working for first to do list item but not for later items.\"\n        break\n    }\n}\n"}
